# Daily attendance processing - 2026-01-17 19:08:44
# Swap the order of names in the "Recorded By" column (G) from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# for every row where that exact text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $v = $cell.Value()
    if ($v -eq $oldText) {
        $cell.Value = $newText
    }
}
